# Rename the inline picture shapes' "name" identifiers.
#
# The document embeds the Pearson Edexcel logo (PearsonLogo.png) twice -
# once in the "first page" footer and once in the default footer - plus
# the BTEC logo (BTec_Logo-Orange) once in the "first page" header.
#
#   - the two PearsonLogo pictures were named "image2.png" and are
#     renamed to "image1.png"
#   - the BTec_Logo-Orange picture was named "image1.jpg" and is
#     renamed to "image2.jpg"
#
# InlineShape objects don't expose a settable Name property (same as in
# real Word), so each picture is temporarily converted to a floating
# Shape - which does expose .Name - renamed, and converted back to an
# inline shape so the on-disk <wp:inline> layout is preserved.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Footer used on the first page (footer1.xml) - PearsonLogo image2.png -> image1.png
Rename-InlinePicture $sec.Footers.Item(2).Range "image1.png"

# Default footer (footer2.xml) - PearsonLogo image2.png -> image1.png
Rename-InlinePicture $sec.Footers.Item(1).Range "image1.png"

# Header used on the first page (header1.xml) - BTec_Logo-Orange image1.jpg -> image2.jpg
Rename-InlinePicture $sec.Headers.Item(2).Range "image2.jpg"
